$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the old closing "meta description" text (the italic paragraph
#    right at the end of the document) with the new DALL-E image prompt.
#    Done first, while the phrase is still unique in the document.
# ---------------------------------------------------------------------------
$oldText = "Read our review of Dragon Egg, an online slot game with a free spins mode. Enjoy higher value wins than average and clear graphics. Play Dragon Egg for free."
$newText = "Prompt for DALLE: Create a feature image for Dragon Egg, a slot game by Tom Horn, in a cartoon style. The image should feature a happy Maya warrior wearing glasses. The warrior should have a confident expression on their face and be holding a golden dragon egg in one hand, as if they have just won it in the slot game. The background should be a dark cave, with shadows of dragons visible in the background. The image should be eye-catching and convey the excitement of winning big in the game."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph that used to sit just before
#    that closing paragraph near the end of the document.
# ---------------------------------------------------------------------------
$dupPara = $null
For ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text.TrimEnd([char]13)
    if ($txt -eq "Play Dragon Egg Free Slot Review | Exciting Free Spins Mode" -and $para.Range.Bold) {
        $dupPara = $para
        break
    }
}
if ($dupPara -ne $null) {
    $dupPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range.Duplicate
$titleRange.Collapse(0)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Dragon Egg Free Slot Review | Exciting Free Spins Mode</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Dragon Egg, an online slot game with a free spins mode. Enjoy higher value wins than average and clear graphics. Play Dragon Egg for free.</w:t></w:r></w:p>'

$titleRange.InsertXML($metaXml)
